$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-91 down to 28-92.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new data record.
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 44497
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 100112012
$ws.Range("G27").Value = "Espinaca"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 7000
$ws.Range("N27").Value = "$/docena de atados"
$ws.Range("O27").Value = "Región de La Araucanía"
$ws.Range("P27").Value = 2333
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = "Hortaliza"
